$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reporting-period dates in row 8 (cambios de agosto / fe de
# erratas, historico): start/end of period, validation and update dates.
$ws.Range("B8").Value = 44652
$ws.Range("C8").Value = 44742
$ws.Range("F8").Value = 44753
$ws.Range("G8").Value = 44753

# Update the sheet view: scroll the window so A2 is the top-left visible
# cell, and move the selection to C12.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12").Select()
